# Applies the "done with client and admin schema api" edit:
#  - Highlights three Client Schema bullet paragraphs yellow
#  - Updates "Update client data by id" -> "Update client data by slug"
#  - Splits "profile" in the User Schema bullet #3 to relocate the
#    "_GoBack" bookmark inside that word
#  - Removes a now-superfluous trailing empty paragraph near the end
#    of the document
$d = $word.ActiveDocument

# --- 1. "3. Client can see only his activity by year. (rest)" -> yellow highlight
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -eq "3. Client can see only his activity by year. (rest)`r") {
        $p.Range.HighlightColorIndex = 7
        break
    }
}

# --- 2. "4. Update client data by id (gql)" -> "4. Update client data by slug (gql)" + yellow highlight
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -eq "4. Update client data by id (gql)`r") {
        $p.Range.Find.Execute("by id ", $true, $false, $false, $false, $false, $true, 1, $false, "by slug ", 2)
        $p.Range.HighlightColorIndex = 7
        break
    }
}

# --- 3. "5. Client can see only his published or draft blog (gql)" -> yellow highlight
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -eq "5. Client can see only his published or draft blog (gql)`r") {
        $p.Range.HighlightColorIndex = 7
        break
    }
}

# --- 4. Relocate the "_GoBack" bookmark into "3. User can see his profile (rest)",
#        splitting right after the "p" of "profile".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -eq "3. User can see his profile (rest)`r") {
        $paraStart = $p.Range.Start
        $splitAt = $paraStart + ("3. User can see his p").Length

        $old = $d.Bookmarks("_GoBack")
        $old.Delete()

        $splitRange = $d.Range($splitAt, $splitAt)
        $d.Bookmarks.Add("_GoBack", $splitRange)
        break
    }
}

# --- 5. Remove the extra empty paragraph right before the final empty
#        paragraph / sectPr at the end of the document.
$n = $d.Paragraphs.Count
$last = $d.Paragraphs($n)
$secondLast = $d.Paragraphs($n - 1)
if ($last.Range.Text -eq "`r" -and $secondLast.Range.Text -eq "`r") {
    $secondLast.Range.Delete()
}
